# Daily attendance processing - 2026-01-09 06:47:29
# Normalises the "Recorded By" column (G) so that any "System" entry is
# moved to the front of the comma-separated recorder list, preserving the
# relative order of the remaining recorders. Cells without a "System"
# entry still need reordering (observed behaviour: the list is produced
# most-recent-first), cells that already start with "System" (or have only
# a single recorder) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Transform-RecordedBy($old) {
    $parts = $old -split ", "
    if ($parts.Length -le 1) {
        return $old
    }

    $sysParts = @()
    $otherParts = @()
    foreach ($p in $parts) {
        if ($p.Equals("System")) {
            $sysParts += $p
        } else {
            $otherParts += $p
        }
    }

    if ($sysParts.Length -gt 0) {
        $newParts = $sysParts + $otherParts
    } else {
        $newParts = @()
        for ($i = $parts.Length - 1; $i -ge 0; $i--) {
            $newParts += $parts[$i]
        }
    }

    return [string]::Join(", ", $newParts)
}

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $old = $cell.Text
    if ([string]::IsNullOrEmpty($old)) {
        continue
    }
    $new = Transform-RecordedBy($old)
    if (-not $new.Equals($old)) {
        $cell.Value = $new
    }
}
